$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.993.97"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.549.71"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "305.18"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "98.15"
$ws.Range("E6").Value = "  +6.68%  "
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "37.01"
$ws.Range("E10").Value = "  +3.35%  "
$ws.Range("D11").Value = "0.0838"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.944.52"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "2.550.77"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "15.18"
$ws.Range("E16").Value = "  +7.91%  "
$ws.Range("D17").Value = "0.877"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "43.015.43"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "13.84"
$ws.Range("E19").Value = "  +6.47%  "
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "72.03"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "255.67"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "2.09"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "28.13"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "10.27"
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("D29").Value = "37.94"
$ws.Range("D30").Value = "6.23"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").Value = "158.71"
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("D33").Value = "19.76"
$ws.Range("E33").Value = "  +17.22%  "
$ws.Range("D34").Value = "2.15"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "0.0805"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "25.56"
$ws.Range("E39").Value = "  +8.23%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("E43").Value = "  +28.17%  "
$ws.Range("D44").Value = "2.104.36"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "86.93"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "2.802.03"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "74.91"
$ws.Range("E50").Value = "  +9.21%  "
$ws.Range("D51").Value = "103.80"
$ws.Range("E51").Value = "  -0.17%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
